$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(8,  "acetylcystein 200 mg f 4"),
    @(9,  "al . £"),
    @(10, "acetylcystein 200 mg a ze"),
    @(11, "acetylcystein 200 mg a"),
    @(12, "‘acehasan 200"),
    @(13, "acehasan, 200. _"),
    @(14, "acetylcystein ay)"),
    @(15, "acetylcystein 200 mg a"),
    @(16, "atp.")
)

$row = 9
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = "Thuốc"
    $row++
}
